$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room for the new rows. The sheet used to end at row 31 (the closing
#    border row). The updated report keeps that closing border row but now at
#    row 41, with freshly-blank rows 31-40 above it (same look as rows 23-30).
# ---------------------------------------------------------------------------

# Push the old "closing" row (31) format down to the new last row (41).
$ws.Range("A31:C31").Copy() | Out-Null
$ws.Range("A41:C41").PasteSpecial(-4122) | Out-Null

# Fill rows 31-40 with the plain blank-row formatting used by rows 23-30.
$ws.Range("A30:C30").Copy() | Out-Null
$ws.Range("A31:C40").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Populate the newly-used rows (12-22) with date / hours / activity data.
# ---------------------------------------------------------------------------

$rows = @(
    @{ Row = 12; Date = 42825; Hours = 2;    Text = "Searched for HTML library for java" },
    @{ Row = 13; Date = 42816; Hours = 2;    Text = "Created classes to represent database objects" },
    @{ Row = 14; Date = 42817; Hours = 2;    Text = "Writing Database Interface Functions & Organized SQL Scripts" },
    @{ Row = 15; Date = 42832; Hours = 0.5;  Text = "Compiled existing SQL scripts into one file" },
    @{ Row = 16; Date = 42834; Hours = 0.5;  Text = "Wrote SQL Script for resetting Database" },
    @{ Row = 17; Date = 42836; Hours = 2;    Text = "Writing Database Interface Functions & Writing JSP Pages" },
    @{ Row = 18; Date = 42841; Hours = 1.5;  Text = "Writing Database Interface Functions" },
    @{ Row = 19; Date = 42842; Hours = 2;    Text = "Create Abstract Page Class" },
    @{ Row = 20; Date = 42843; Hours = 1.5;  Text = "Writing Database Interface Functions" },
    @{ Row = 21; Date = 42848; Hours = 0.5;  Text = "Reviewing Other Group Member's Code" },
    @{ Row = 22; Date = 42852; Hours = 14;   Text = "Finishing Code" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.Hours
}

# Row 22's date keeps the "m/d/yy" look already used higher up the sheet
# (same style as row 3, etc.), while rows 12-21 use the "d-mmm" look already
# used by row 11.
$ws.Cells.Item(22, 1).NumberFormat = "m/d/yy"
foreach ($rowNum in 12..21) {
    $ws.Cells.Item($rowNum, 1).NumberFormat = "d-mmm"
}

# Activity text - set in this specific order so identical strings merge into
# the same shared-string entry the same way the original report grew.
$ws.Cells.Item(22, 3).Value = "Finishing Code"
$ws.Cells.Item(21, 3).Value = "Reviewing Other Group Member's Code"
$ws.Cells.Item(18, 3).Value = "Writing Database Interface Functions"
$ws.Cells.Item(19, 3).Value = "Create Abstract Page Class"
$ws.Cells.Item(17, 3).Value = "Writing Database Interface Functions & Writing JSP Pages"
$ws.Cells.Item(16, 3).Value = "Wrote SQL Script for resetting Database"
$ws.Cells.Item(15, 3).Value = "Compiled existing SQL scripts into one file"
$ws.Cells.Item(12, 3).Value = "Searched for HTML library for java"
$ws.Cells.Item(14, 3).Value = "Writing Database Interface Functions & Organized SQL Scripts"
$ws.Cells.Item(13, 3).Value = "Created classes to represent database objects"
$ws.Cells.Item(20, 3).Value = "Writing Database Interface Functions"

# ---------------------------------------------------------------------------
# 3) Restore the on-screen selection to where the author left off (C14).
# ---------------------------------------------------------------------------
$ws.Range("C14").Select() | Out-Null
